$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "years" table (row 4 header, rows 5-8 data) is being extended one
# column to the right: a new column Q is inserted right after the last
# year column (P/2019) to hold the 2020 figures. Inserting the column
# (rather than just typing into the blank cell) picks up column P's
# per-row styling (number format / font / borders) for the new cells,
# matching the existing "year" columns.
$ws.Range("Q1").EntireColumn.Insert()

# Fill in the 2020 values - same figures as reported for 2019, aside from
# the year header itself.
$ws.Range("Q4").Value = 2020
$ws.Range("Q5").Value = 2
$ws.Range("Q6").Value = 0.3
$ws.Range("Q7").Value = 0.1
$ws.Range("Q8").Value = 4.3

# Leave the selection where the author ended up after editing.
$ws.Range("O12").Select()
